$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the color scale ("cols") values for the tasmin/tasmax rows
# from "red, yellow, blue" to "blue, yellow, red" (issue #43: fixing color scale)
$ws.Range("D5").Value = "blue, yellow, red"
$ws.Range("D6").Value = "blue, yellow, red"

# Update the active selection to D6 as reflected in the saved file
$ws.Activate()
$ws.Range("D6").Select()
